$d = $word.ActiveDocument

# --- 1) "season tickets, etc..." run clean-up -----------------------------
# Word's spell-checker had split "season tickets, etc" across several runs
# (with proofErr spell-check markers around "etc"). Re-typing/correcting
# that phrase collapses it back into a single run and drops the proofErr
# markers.
$d.Content.Find.Execute("season tickets, etc…) in the preferences setup are not considered.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "season tickets, etc…) in the preferences setup are not considered.", 2) | Out-Null

# --- 2) BCE diagrams paragraph text edits ---------------------------------
# "...is adopted. BCE diagrams..." -> "...is adopted and BCE diagrams..."
$d.Content.Find.Execute("is adopted. BCE diagrams", $true, $false, $false, $false, $false,
                         $true, 1, $false, "is adopted and BCE diagrams", 2) | Out-Null

# Add "object" after "Entities" and after "Controls"
$d.Content.Find.Execute("Entities model the access to data; Controls manage the communication", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Entities object model the access to data; Controls object manage the communication", 2) | Out-Null

# --- 3) Move the _GoBack bookmark -----------------------------------------
# In the original the _GoBack bookmark (empty range) sat between "...the
# system" and ". Boundaries...". After the edit above it should wrap the
# whole edited sentence, i.e. start right before "For the implementation"
# and end right after "...entities." (mirroring where Word leaves _GoBack
# after the last edit in the paragraph touched both ends of the sentence).
$startRng = $d.Content
$startRng.Find.Execute("For the implementation", $true, $false, $false, $false, $false,
                        $true, 1, $false, "", 0) | Out-Null
$startPos = $startRng.Start

$endRng = $d.Content
$endRng.Find.Execute("manage the communication between boundaries and entities.", $true, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$endPos = $endRng.End

$bookmarkRange = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
